# Bug with translation of "&"
#
# The lookup key "Initial & Final Surveillance Diagnosis" (used verbatim,
# literal ampersand) did not match what the app was actually searching for,
# so the translation never resolved. Fix: add a new translation row that
# spells the key out with "and" instead of "&" (and keep the already-correct
# French translation for it). The old "&"-keyed row is kept (for backward
# compatibility / until it is cleaned up) but its French cell is reset to a
# "TBT" (to-be-translated) placeholder since it's no longer the canonical
# entry.
#
# Same root cause affected "Susceptible & Intermediate are always combined
# in this visualisation of co-resistances." - fixed in place by spelling out
# "and"/"et" instead of using "&".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 75, pushing the previous row 75 ("Isolates"/"Isolats")
# and everything below it down by one (old 75..181 -> new 76..182).
$ws.Rows("75:75").Insert()

# New row 75: the corrected lookup key (spelled-out "and") paired with the
# French text that used to live on row 74.
$ws.Range("A75").Value = "Initial and Final Surveillance Diagnosis"
$ws.Range("B75").Value = "Diagnostic de surveillance initial et final"

# Row 74 keeps its original (ampersand) English key but its French
# translation becomes a placeholder pending cleanup.
$ws.Range("B74").Value = "TBT"

# Fix the other "&" translation bug in place (row shifted down to 145 by
# the insert above).
$found = $ws.Cells.Find("Susceptible & Intermediate are always combined in this visualisation of co-resistances.")
if ($found -ne $null) {
    $r = $found.Row
    $ws.Cells.Item($r, 1).Value = "Susceptible and Intermediate are always combined in this visualisation of co-resistances."
    $ws.Cells.Item($r, 2).Value = "Susceptible et Intermédiaire sont toujours combinés dans cette visualisation des co-résistances."
}
